$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content but keep formatting in place (A1:E1 keep their header style)
$ws.Cells.ClearContents()

# Header row: insert "accuracy" right after "model" (before sensitivity/specificity/precision/smote)
$ws.Cells.Item(1, 1).Value = "model"
$ws.Cells.Item(1, 2).Value = "accuracy"
$ws.Cells.Item(1, 3).Value = "sensitivity"
$ws.Cells.Item(1, 4).Value = "specificity"
$ws.Cells.Item(1, 5).Value = "precision"
# F1 needs the same header styling as A1:E1 - copy format from E1 before assigning its value
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Cells.Item(1, 6).Value = "smote"

# Row 2: gradient_boosting
$ws.Cells.Item(2, 1).Value = "gradient_boosting"
$ws.Cells.Item(2, 2).Value = 0.9970788704965921
$ws.Cells.Item(2, 3).Value = 0.9655172413793104
$ws.Cells.Item(2, 4).Value = 0.9979959919839679
$ws.Cells.Item(2, 5).Value = 0.9333333333333333
$ws.Cells.Item(2, 6).Value = $true

# Row 3: logistic_regression
$ws.Cells.Item(3, 1).Value = "logistic_regression"
$ws.Cells.Item(3, 2).Value = 0.9970788704965921
$ws.Cells.Item(3, 3).Value = 0.896551724137931
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = $true
